$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14) - existing N/O/P shift right to O/P/Q.
$ws.Columns.Item(14).Insert()

# New column inherits the neighbouring (M) column's width.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Update the sheet's selection and make it the active/selected sheet & tab.
$ws.Range("L13").Select() | Out-Null
$ws.Activate()
